# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" on all three sheets
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. Narrow the now-shorter "Status"/language columns:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status values -------------------------------------------------
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value     = "In Translation"
$wsDeDe.Range("C2:C4").Value     = "In Translation"

# --- Shrink the affected columns ------------------------------------------
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newWidth
